# Applies the "Add files via upload" change: appends a new results table
# (rows 17-23) to Hoja1, tweaks a couple of existing cells/column widths,
# and updates the saved selections on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- existing cell I14 now holds a value -----------------------------------
$ws1.Cells.Item(14, 9).Value = 1

# --- new header row (row 17) ------------------------------------------------
$ws1.Cells.Item(17, 2).Value = "Dataset"
$ws1.Cells.Item(17, 3).Value = "Metodo"
$ws1.Cells.Item(17, 4).Value = "Distancia"
$ws1.Cells.Item(17, 5).Value = "Indice Cofenetico"
$ws1.Cells.Item(17, 6).Value = "Indice Rand Dist Cant Clusters "
$ws1.Cells.Item(17, 7).Value = "Indice Rand Distance"
$ws1.Cells.Item(17, 8).Value = "Rand MaxClust 9"

# --- new data rows (18-22) --------------------------------------------------
$ws1.Cells.Item(18, 1).Value = "DF1"
$ws1.Cells.Item(18, 2).Value = "Todas la variables"
$ws1.Cells.Item(18, 3).Value = "average"
$ws1.Cells.Item(18, 4).Value = "euclidea"
$ws1.Cells.Item(18, 5).Value = 0.625667903586
$ws1.Cells.Item(18, 6).Value = 8
$ws1.Cells.Item(18, 7).Value = 0.00050339369652909003
$ws1.Cells.Item(18, 8).Value = 0.00058923346526635104

$ws1.Cells.Item(19, 1).Value = "DF2"
$ws1.Cells.Item(19, 2).Value = "Timbres y Pitches"
$ws1.Cells.Item(19, 3).Value = "average"
$ws1.Cells.Item(19, 4).Value = "euclidea"
$ws1.Cells.Item(19, 5).Value = 0.67914832900698097
$ws1.Cells.Item(19, 6).Value = 7
$ws1.Cells.Item(19, 7).Value = 0.00058496593847128799
$ws1.Cells.Item(19, 8).Value = 0.00116786669924672

$ws1.Cells.Item(20, 1).Value = "DF3"
$ws1.Cells.Item(20, 2).Value = "Timbres  "
$ws1.Cells.Item(20, 3).Value = "average"
$ws1.Cells.Item(20, 4).Value = "euclidea"
$ws1.Cells.Item(20, 5).Value = 0.62696867077538998
$ws1.Cells.Item(20, 6).Value = 7
$ws1.Cells.Item(20, 7).Value = 0.000045421215857999997
$ws1.Cells.Item(20, 8).Value = 0.000193636861092818

$ws1.Cells.Item(21, 1).Value = "DF4"
$ws1.Cells.Item(21, 2).Value = "Pitches"
$ws1.Cells.Item(21, 3).Value = "average"
$ws1.Cells.Item(21, 4).Value = "euclidea"
$ws1.Cells.Item(21, 5).Value = 0.69966714108047501
$ws1.Cells.Item(21, 6).Value = 6
$ws1.Cells.Item(21, 7).Value = 0.00032227911394587999
$ws1.Cells.Item(21, 8).Value = 0.00046595291814487399

$ws1.Cells.Item(22, 1).Value = "DF5"
$ws1.Cells.Item(22, 2).Value = "Audio Features"
$ws1.Cells.Item(22, 3).Value = "average"
$ws1.Cells.Item(22, 4).Value = "euclidea"
$ws1.Cells.Item(22, 5).Value = 0.626671566854476
$ws1.Cells.Item(22, 6).Value = 5
$ws1.Cells.Item(22, 7).Value = 0.00046039562690700097
$ws1.Cells.Item(22, 8).Value = 0.00043567739411282

# --- distinct cell formats reused from the rest of the sheet ---------------
# style "s=1" (font color 212121, no scheme)
$ws1.Cells.Item(14, 8).Copy() | Out-Null
$ws1.Cells.Item(19, 5).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(19, 8).PasteSpecial(-4122) | Out-Null

# style "s=3" (font color 212121, minor scheme)
$ws1.Cells.Item(10, 7).Copy() | Out-Null
$ws1.Cells.Item(21, 5).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(21, 7).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(21, 8).PasteSpecial(-4122) | Out-Null

# style "s=2" (used on most filler/blank cells in this table)
$ws1.Cells.Item(9, 2).Copy() | Out-Null
$ws1.Cells.Item(22, 9).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 2).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 3).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 4).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 5).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 6).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 7).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 8).PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(23, 9).PasteSpecial(-4122) | Out-Null

# new number format (0.000000000) combined with the 212121/minor-scheme font
$ws1.Cells.Item(20, 7).Font.Color = 2171169
$ws1.Cells.Item(20, 7).NumberFormat = "0.000000000"

# --- column widths -----------------------------------------------------------
$ws1.Columns.Item(5).ColumnWidth = 27.166666666666664
$ws1.Columns.Item(6).ColumnWidth = 23.833333333333336
$ws1.Columns.Item(7).ColumnWidth = 22.0

# --- selections / view ------------------------------------------------------
$ws2.Range("I18").Select()
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws1.Range("F19").Select()
